# Updates the cryptos list on the active worksheet with latest scraped
# prices / volume percentages, mirroring the GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string cell setter (safe for values Excel won't reinterpret
# as a number, e.g. multi-dot prices, percentages, names, URLs).
function Set-Cell {
    param(
        [string]$Address,
        [string]$Value
    )
    $ws.Range($Address).Value = $Value
}

# Text cell setter for values that "look" numeric (single decimal point)
# so Excel keeps storing them as text instead of converting to a number,
# while not leaving a lingering custom cell style behind.
function Set-TextCell {
    param(
        [string]$Address,
        [string]$Value
    )
    $ws.Range($Address).NumberFormat = "@"
    $ws.Range($Address).Value = $Value
    $ws.Range($Address).Style = "Normal"
}

# Row 2 - Bitcoin
Set-Cell "D2" "26.331.44"
Set-Cell "E2" "  +1.41%  "

# Row 3 - Ethereum
Set-Cell "D3" "1.623.88"
Set-Cell "E3" "  +1.66%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  -0.09%  "

# Row 5 - BNB
Set-TextCell "D5" "212.63"
Set-Cell "E5" "  +0.79%  "

# Row 6 - USDC
Set-Cell "E6" "  -0.09%  "

# Row 7 - XRP
Set-TextCell "D7" "0.488"
Set-Cell "E7" "  +1.15%  "

# Row 8 - Cardano
Set-TextCell "D8" "0.249"
Set-Cell "E8" "  +1.41%  "

# Row 9 - Dogecoin
Set-TextCell "D9" "0.0615"
Set-Cell "E9" "  +0.73%  "

# Row 10 - Solana
Set-TextCell "D10" "18.96"
Set-Cell "E10" "  +5.09%  "

# Row 11 - TRON
Set-TextCell "D11" "0.0816"
Set-Cell "E11" "  +0.63%  "

# Row 12 and 13 swap places: WrappedliquidstakedEther2.0 <-> WrappedEther
Set-Cell "B12" "WrappedEther"
Set-Cell "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-Cell "D12" "1.675.68"
Set-Cell "E12" "  +4.87%  "

Set-Cell "B13" "WrappedliquidstakedEther2.0"
Set-Cell "C13" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-Cell "D13" "1.850.65"
Set-Cell "E13" "  +1.63%  "

# Row 14 - Polkadot
Set-Cell "E14" "  +0.87%  "

# Row 15 - Polygon
Set-Cell "E15" "  +1.27%  "

# Row 16 - WrappedBTC
Set-Cell "D16" "26.339.25"
Set-Cell "E16" "  +1.39%  "

# Row 17 - Litecoin
Set-Cell "E17" "  +4.01%  "

# Row 18 - ShibaInu
Set-Cell "E18" "  +1.12%  "

# Row 19 - Dai
Set-Cell "E19" "  -0.10%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "203.05"
Set-Cell "E20" "  +1.41%  "

# Row 21 - Uniswap
Set-Cell "E21" "  +1.80%  "

# Row 22 - Avalanche
Set-Cell "E22" "  +1.36%  "

# Row 23 - Chainlink
Set-TextCell "D23" "6.07"
Set-Cell "E23" "  +1.14%  "

# Row 24 - Toncoin
Set-Cell "E24" "  +7.45%  "

# Row 25 - Monero
Set-TextCell "D25" "143.24"
Set-Cell "E25" "  +0.92%  "

# Row 26 - BinanceUSD
Set-Cell "E26" "  -0.10%  "

# Row 27 - Stellar
Set-Cell "E27" "  +0.23%  "

# Row 28 - EthereumClassic
Set-Cell "E28" "  +0.91%  "

# Row 29 - Cosmos
Set-TextCell "D29" "6.58"
Set-Cell "E29" "  +1.85%  "

# Row 30 - Hedera
Set-TextCell "D30" "0.0528"
Set-Cell "E30" "  +11.02%  "

# Row 31 - PancakeSwap
Set-Cell "E31" "  +0.85%  "

# Row 32 - Filecoin
Set-Cell "E32" "  +2.91%  "

# Row 33 - InternetComputer(DFINITY)
Set-Cell "E33" "  +0.10%  "

# Row 34 and 35 swap places: LidoDAOToken <-> HuobiToken
Set-Cell "B34" "HuobiToken"
Set-Cell "C34" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D34" "2.44"
Set-Cell "E34" "  +3.24%  "

Set-Cell "B35" "LidoDAOToken"
Set-Cell "C35" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D35" "1.51"
Set-Cell "E35" "  +1.98%  "

# Row 36 - Maker
Set-Cell "D36" "1.176.81"
Set-Cell "E36" "  +4.74%  "

# Row 37 - VeChain
Set-Cell "E37" "  +1.39%  "

# Row 38 - ARBITRUM
Set-TextCell "D38" "0.815"
Set-Cell "E38" "  +3.61%  "

# Row 39 - PaxDollar
Set-Cell "E39" "  -0.06%  "

# Row 40 - MXToken
Set-Cell "E40" "  +0.08%  "

# Row 41 - ImmutableX
Set-TextCell "D41" "0.499"
Set-Cell "E41" "  +1.79%  "

# Row 42 - TrustWalletToken
Set-Cell "E42" "  +1.36%  "

# Row 43 - FraxShare
Set-Cell "E43" "  +4.28%  "

# Row 44 - RocketPoolETH
Set-Cell "D44" "1.761.16"
Set-Cell "E44" "  +1.66%  "

# Row 45 - Quant
Set-TextCell "D45" "93.47"
Set-Cell "E45" "  +0.76%  "

# Row 46 - BabyDogeCoin
Set-Cell "E46" "  +15.43%  "

# Row 47 - RenderToken
Set-Cell "E47" "  +1.43%  "

# Row 48 - Aave
Set-TextCell "D48" "54.12"
Set-Cell "E48" "  +1.24%  "

# Row 49 - Cronos
Set-Cell "E49" "  +0.98%  "

# Row 50 - Mantle
Set-Cell "E50" "  +0.03%  "

# Row 51 - USDD
Set-Cell "E51" "  -0.45%  "
